$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Token Usage"
$ws.Name = "Token Usage"

# Clear the bold/border formatting previously applied to A2 ("Name" header cell)
$ws.Range("A2").ClearFormats()

# Row 1: column headers
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "GPT4o"
$ws.Range("C1").Value = "Gemini"
$ws.Range("D1").Value = "LLaMA"
$ws.Range("E1").Value = "Claude"

# Row 2: Conveyor_System_Controller token usage
$ws.Range("A2").Value = "Conveyor_System_Controller"
$ws.Range("B2").Value = 979
$ws.Range("C2").Value = 9589
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1320

# Row 3: Traffic_Light_Controller token usage
$ws.Range("A3").Value = "Traffic_Light_Controller"
$ws.Range("B3").Value = 804
$ws.Range("C3").Value = 6819
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1073
